$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.813.60'
$ws.Range('E2').Value = '  -2.00%  '
$ws.Range('D3').Value = '2.271.44'
$ws.Range('E3').Value = '  -3.30%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = '''309.14'
$ws.Range('E5').Value = '  -4.90%  '
$ws.Range('D6').Value = '''104.70'
$ws.Range('E6').Value = '  +3.12%  '
$ws.Range('D7').Value = '''0.623'
$ws.Range('E7').Value = '  -2.21%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').Value = '''0.600'
$ws.Range('E9').Value = '  -3.32%  '
$ws.Range('D10').Value = '''40.07'
$ws.Range('E10').Value = '  +0.27%  '
$ws.Range('D11').Value = '''0.0903'
$ws.Range('E11').Value = '  -1.97%  '
$ws.Range('D12').Value = '''8.19'
$ws.Range('E12').Value = '  -2.70%  '
$ws.Range('E13').Value = '  -0.19%  '
$ws.Range('D14').Value = '''0.956'
$ws.Range('E14').Value = '  -3.83%  '
$ws.Range('D15').Value = '''15.34'
$ws.Range('E15').Value = '  -4.63%  '
$ws.Range('D16').Value = '2.613.76'
$ws.Range('E16').Value = '  -3.40%  '
$ws.Range('D17').Value = '2.270.04'
$ws.Range('E17').Value = '  -3.38%  '
$ws.Range('D18').Value = '41.776.79'
$ws.Range('E18').Value = '  -1.97%  '
$ws.Range('D19').Value = '''7.50'
$ws.Range('E19').Value = '  -5.21%  '
$ws.Range('E20').Value = '  -2.58%  '
$ws.Range('D21').Value = '''73.18'
$ws.Range('E21').Value = '  -4.28%  '
$ws.Range('D22').Value = '''3.41'
$ws.Range('E22').Value = '  -8.08%  '
$ws.Range('D23').Value = '''254.48'
$ws.Range('E23').Value = '  -3.59%  '
$ws.Range('D24').Value = '''2.28'
$ws.Range('E24').Value = '  -1.09%  '
$ws.Range('D25').Value = '''9.20'
$ws.Range('E25').Value = '  -8.61%  '
$ws.Range('E26').Value = '  +0.55%  '
$ws.Range('D27').Value = '''10.86'
$ws.Range('E27').Value = '  -5.06%  '
$ws.Range('E28').Value = '  +3.23%  '
$ws.Range('D29').Value = '''22.49'
$ws.Range('E29').Value = '  -0.32%  '
$ws.Range('D30').Value = '''165.40'
$ws.Range('E30').Value = '  -5.76%  '
$ws.Range('D31').Value = '''35.21'
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('D32').Value = '''0.0882'
$ws.Range('E32').Value = '  -1.63%  '
$ws.Range('D33').Value = '''2.92'
$ws.Range('E33').Value = '  -5.93%  '
$ws.Range('D34').Value = '''5.72'
$ws.Range('E34').Value = '  -5.46%  '
$ws.Range('D35').Value = '''0.128'
$ws.Range('E35').Value = '  -2.52%  '
$ws.Range('D36').Value = '''0.116'
$ws.Range('E36').Value = '  +6.78%  '
$ws.Range('D37').Value = '''4.51'
$ws.Range('E37').Value = '  -1.13%  '
$ws.Range('D38').Value = '''0.0349'
$ws.Range('E38').Value = '  -2.48%  '
$ws.Range('D39').Value = '''2.72'
$ws.Range('E39').Value = '  -3.36%  '
$ws.Range('D40').Value = '''3.59'
$ws.Range('E40').Value = '  -4.72%  '
$ws.Range('D41').Value = '''71.21'
$ws.Range('E41').Value = '  +2.05%  '
$ws.Range('D42').Value = '''96.68'
$ws.Range('E42').Value = '  +4.27%  '
$ws.Range('D43').Value = '''1.44'
$ws.Range('E43').Value = '  -3.14%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').Value = '''1.00'
$ws.Range('E44').Value = '  +0.32%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').Value = '''0.225'
$ws.Range('E45').Value = '  -5.18%  '
$ws.Range('D46').Value = '''12.13'
$ws.Range('E46').Value = '  +2.31%  '
$ws.Range('D47').Value = '''110.81'
$ws.Range('E47').Value = '  -8.59%  '
$ws.Range('D48').Value = '''8.94'
$ws.Range('E48').Value = '  -2.49%  '
$ws.Range('D49').Value = '''5.24'
$ws.Range('E49').Value = '  -4.98%  '
$ws.Range('D50').Value = '''73.86'
$ws.Range('E50').Value = '  +4.60%  '
$ws.Range('D51').Value = '1.552.59'
$ws.Range('E51').Value = '  +0.24%  '
